# Tour de France tracking workbook update:
#  - "Uitvallers" (riders who dropped out) gains 4 new names
#  - "Huidig" (current standings) gets column O filled in for the new stage

$wb = $excel.ActiveWorkbook

# --- Sheet "Uitvallers": append riders who abandoned the race ---
$uitvallers = $wb.Worksheets.Item("Uitvallers")
$uitvallers.Range("A15").Value = "Remco Evenepoel"
$uitvallers.Range("A16").Value = "Mattias Skjelmose"
$uitvallers.Range("A17").Value = "Bryan Coquard"
$uitvallers.Range("A18").Value = "Steff Cras"

# --- Sheet "Huidig": fill column O (rows 6-19) with this stage's picks ---
$huidig = $wb.Worksheets.Item("Huidig")
$huidig.Range("O6").Value = "Thymen Arensman"
$huidig.Range("O7").Value = "Tadej Pogacar"
$huidig.Range("O8").Value = "Jonas Vingegaard"
$huidig.Range("O9").Value = "Félix Gall"
$huidig.Range("O10").Value = "Florian Lipowitz"
$huidig.Range("O11").Value = "Oscar Onley"
$huidig.Range("O12").Value = "Ben Healy"
$huidig.Range("O13").Value = "Primoz Roglic"
$huidig.Range("O14").Value = "Tobias Johannessen"
$huidig.Range("O15").Value = "Kévin Vauquelin"
$huidig.Range("O16").Value = "Tadej Pogacar"
$huidig.Range("O17").Value = "Jonathan Milan"
$huidig.Range("O18").Value = "Lenny Martinez"
$huidig.Range("O19").Value = "Florian Lipowitz"

# Leave the view focused on "Huidig", with O19 as the active cell
# (the new cell the user last filled in), and drop the selection on
# "Uitvallers" to the row right after the new data.
$uitvallers.Range("A19").Select() | Out-Null
$huidig.Activate() | Out-Null
$huidig.Range("O19").Select() | Out-Null
